$wb = $excel.ActiveWorkbook

# --- "Project Ideas" sheet (sheet2) first, so new shared strings land in the
# --- same order the original authoring tool produced them ---
$ideas = $wb.Worksheets.Item("Project Ideas")
$ideas.Range("A11").Value = "ITC over multiple years"

# --- "To Do" sheet (sheet1) ---
$todo = $wb.Worksheets.Item("To Do")

# Status column (A) updates: "Not done" -> "Done" for several rows now finished
$todo.Range("A36").Value = "Done"
$todo.Range("A38").Value = "Done"
$todo.Range("A39").Value = "Done"

# Status column (A) updates: "Not done" -> "POUT" (new status value)
$todo.Range("A45").Value = "POUT"
$todo.Range("A46").Value = "POUT"

# Status column (A) updates: "Not done" -> "Future"
$todo.Range("A49").Value = "Future"

# Fill in "Who" column (C) for a couple of rows that were blank
$todo.Range("C52").Value = "Steve"
$todo.Range("C60").Value = "Janine"

# More "Not done" -> "Done" updates
$todo.Range("A58").Value = "Done"
$todo.Range("A63").Value = "Done"
$todo.Range("A64").Value = "Done"

# More "Not done" -> "Future" updates
$todo.Range("A65").Value = "Future"

$todo.Range("A66").Value = "Done"

$todo.Range("A67").Value = "Future"
$todo.Range("A68").Value = "Future"
$todo.Range("A74").Value = "Future"

$todo.Range("A111").Value = "Done"
$todo.Range("A112").Value = "Done"

# New status entries for rows 123-125 (meeting notes added at bottom of list)
$todo.Range("A123").Value = "Done"
$todo.Range("A124").Value = "Not done"
$todo.Range("A125").Value = "Not done"

$ideas.Range("A12").Select()

# --- view / selection state ---
# Make "To Do" the active sheet/tab (matches activeTab defaulting back to sheet 1)
$todo.Activate()
$aw = $excel.ActiveWindow
$aw.FreezePanes = $true
$todo.Range("A125").Select()
